$wb = $excel.ActiveWorkbook
Write-Output "noop"
